$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.304.44"
$ws.Range("E2").Value = "  +1.59%  "
$ws.Range("D3").Value = "2.007.43"
$ws.Range("E3").Value = "  +5.55%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'244.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.28%  "
$ws.Range("D6").Value = "'0.659"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.02%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'44.36"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.44%  "
$ws.Range("D9").Value = "'61.15"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.34%  "
$ws.Range("E10").Value = "  +0.86%  "
$ws.Range("D11").Value = "'0.0712"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.09%  "
$ws.Range("E12").Value = "  -0.82%  "
$ws.Range("D13").Value = "'14.26"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.75%  "
$ws.Range("D14").Value = "2.293.52"
$ws.Range("E14").Value = "  +5.50%  "
$ws.Range("D15").Value = "'0.800"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("D16").Value = "2.002.32"
$ws.Range("E16").Value = "  +5.42%  "
$ws.Range("D17").Value = "'4.86"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.84%  "
$ws.Range("D18").Value = "36.287.69"
$ws.Range("E18").Value = "  +1.78%  "
$ws.Range("D19").Value = "'70.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.99%  "
$ws.Range("D20").Value = "0.0₃0809"
$ws.Range("E20").Value = "  -3.01%  "
$ws.Range("D21").Value = "'236.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.60%  "
$ws.Range("D22").Value = "'12.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.83%  "
$ws.Range("D23").Value = "'4.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.91%  "
$ws.Range("D24").Value = "'1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").Value = "'2.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -9.76%  "
$ws.Range("D26").Value = "'165.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.68%  "
$ws.Range("D27").Value = "'8.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.89%  "
$ws.Range("D28").Value = "'19.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.16%  "
$ws.Range("E29").Value = "  -10.75%  "
$ws.Range("D30").Value = "'0.121"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.06%  "
$ws.Range("D31").Value = "'21.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +45.67%  "
$ws.Range("E32").Value = "  -1.87%  "
$ws.Range("E33").Value = "  -4.26%  "
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("D35").Value = "'1.86"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.04%  "
$ws.Range("D36").Value = "'0.0857"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +16.61%  "
$ws.Range("D37").Value = "'3.96"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.30%  "
$ws.Range("D38").Value = "'2.10"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.14%  "
$ws.Range("D39").Value = "'0.843"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.86%  "
$ws.Range("E40").Value = "  -11.58%  "
$ws.Range("E41").Value = "  -5.91%  "
$ws.Range("D42").Value = "'95.10"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.54%  "
$ws.Range("D43").Value = "'1.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.44%  "
$ws.Range("E44").Value = "  +15.08%  "
$ws.Range("D45").Value = "'15.79"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.39%  "
$ws.Range("D46").Value = "1.300.16"
$ws.Range("E46").Value = "  -1.91%  "
$ws.Range("D47").Value = "'0.0811"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("D48").Value = "'2.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.19%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.186.23"
$ws.Range("E49").Value = "  +5.35%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "'2.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -8.56%  "
$ws.Range("D51").Value = "'3.80"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +14.15%  "
